{"js": "const replacements = [\n  [\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\"],\n  [\"72\u00d726=\", \"61\u00d750=\"],\n  [\"64\u00d724=\", \"77\u00d771=\"],\n  [\"34\u00d790=\", \"46\u00d764=\"],\n  [\"13\u00d790=\", \"23\u00d784=\"],\n  [\"32\u00d728=\", \"30\u00d763=\"],\n  [\"11\u00d789=\", \"17\u00d718=\"],\n  [\"64\u00d727=\", \"30\u00d741=\"],\n  [\"48\u00d750=\", \"84\u00d730=\"],\n  [\"11\u00d782=\", \"62\u00d739=\"],\n  [\"89\u00d723=\", \"55\u00d759=\"],\n  [\"36\u00d742=\", \"80\u00d722=\"],\n  [\"99\u00d755=\", \"59\u00d728=\"],\n  [\"50\u00d735=\", \"69\u00d724=\"],\n  [\"60\u00d780=\", \"70\u00d780=\"],\n  [\"98\u00d717=\", \"17\u00d795=\"],\n  [\"40\u00d728=\", \"18\u00d778=\"],\n  [\"17\u00d742=\", \"17\u00d763=\"],\n  [\"11\u00d751=\", \"27\u00d717=\"],\n  [\"23\u00d722=\", \"72\u00d748=\"],\n  [\"83\u00d723=\", \"76\u00d798=\"],\n  [\"30\u00d727=\", \"91\u00d737=\"],\n  [\"71\u00d720=\", \"81\u00d750=\"],\n  [\"67\u00d751=\", \"77\u00d755=\"],\n  [\"75\u00d755=\", \"36\u00d776=\"],\n  [\"25\u00d745=\", \"72\u00d724=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-02-07 Wednesday\", \"2024-02-08 Thursday\")\n    ,@(\"72\u00d726=\", \"61\u00d750=\")\n    ,@(\"64\u00d724=\", \"77\u00d771=\")\n    ,@(\"34\u00d790=\", \"46\u00d764=\")\n    ,@(\"13\u00d790=\", \"23\u00d784=\")\n    ,@(\"32\u00d728=\", \"30\u00d763=\")\n    ,@(\"11\u00d789=\", \"17\u00d718=\")\n    ,@(\"64\u00d727=\", \"30\u00d741=\")\n    ,@(\"48\u00d750=\", \"84\u00d730=\")\n    ,@(\"11\u00d782=\", \"62\u00d739=\")\n    ,@(\"89\u00d723=\", \"55\u00d759=\")\n    ,@(\"36\u00d742=\", \"80\u00d722=\")\n    ,@(\"99\u00d755=\", \"59\u00d728=\")\n    ,@(\"50\u00d735=\", \"69\u00d724=\")\n    ,@(\"60\u00d780=\", \"70\u00d780=\")\n    ,@(\"98\u00d717=\", \"17\u00d795=\")\n    ,@(\"40\u00d728=\", \"18\u00d778=\")\n    ,@(\"17\u00d742=\", \"17\u00d763=\")\n    ,@(\"11\u00d751=\", \"27\u00d717=\")\n    ,@(\"23\u00d722=\", \"72\u00d748=\")\n    ,@(\"83\u00d723=\", \"76\u00d798=\")\n    ,@(\"30\u00d727=\", \"91\u00d737=\")\n    ,@(\"71\u00d720=\", \"81\u00d750=\")\n    ,@(\"67\u00d751=\", \"77\u00d755=\")\n    ,@(\"75\u00d755=\", \"36\u00d776=\")\n    ,@(\"25\u00d745=\", \"72\u00d724=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
